$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1140
$ws.Range("I98").Value = 1016.4545
$ws.Range("K98").Value = 1016.4545
$ws.Range("M98").Value = 481.5454999999999
$ws.Range("H122").Value = 1140
$ws.Range("I122").Value = 1016.4545
$ws.Range("K122").Value = 3049.3635
$ws.Range("M122").Value = -599.3635000000004
$ws.Range("H132").Value = 2440.9656
$ws.Range("I132").Value = 1356.12
$ws.Range("K132").Value = 4068.36
$ws.Range("M132").Value = -1538.36
$ws.Range("H137").Value = 4079.0789
$ws.Range("I137").Value = 1454.7273
$ws.Range("K137").Value = 4364.1819
$ws.Range("M137").Value = -1814.1819
$ws.Range("H138").Value = 2455
$ws.Range("J138").Value = 3999.1428
$ws.Range("L138").Value = 11997.4284
$ws.Range("N138").Value = -22277.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28396.363
$ws.Range("I32").Value = 15434.886
$ws.Range("K32").Value = 15434.886
$ws.Range("M32").Value = -15147.886
$ws.Range("H45").Value = 564163.25
$ws.Range("I45").Value = 1265835.1
$ws.Range("K45").Value = 1265835.1
$ws.Range("M45").Value = -1265458.1
$ws.Range("H61").Value = 1125.2727
$ws.Range("I61").Value = 850.3684
$ws.Range("K61").Value = 850.3684
$ws.Range("M61").Value = -638.3684
$ws.Range("H74").Value = 1829.2
$ws.Range("I74").Value = 1516.8823
$ws.Range("K74").Value = 1516.8823
$ws.Range("M74").Value = -642.8823
$ws.Range("H77").Value = 1829.2
$ws.Range("I77").Value = 1516.8823
$ws.Range("K77").Value = 7584.4115
$ws.Range("M77").Value = -3216.4115
$ws.Range("H110").Value = 2284.1516
$ws.Range("I110").Value = 2330.125
$ws.Range("J110").Value = 813
$ws.Range("K110").Value = 2330.125
$ws.Range("L110").Value = 813
$ws.Range("M110").Value = -285.125
$ws.Range("N110").Value = -4903
$ws.Range("H122").Value = 1710.3334
$ws.Range("I122").Value = 1567.55
$ws.Range("K122").Value = 4702.65
$ws.Range("M122").Value = -2252.65
$ws.Range("H132").Value = 1748.9459
$ws.Range("I132").Value = 1354.069
$ws.Range("K132").Value = 4062.207
$ws.Range("M132").Value = -1532.207
$ws.Range("H136").Value = 1125.2727
$ws.Range("I136").Value = 850.3684
$ws.Range("K136").Value = 2551.1052
$ws.Range("M136").Value = -1.105199999999968
$ws.Range("H140").Value = 120428.5
$ws.Range("J140").Value = 120428.5
$ws.Range("L140").Value = 120428.5
$ws.Range("N140").Value = -130788.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2162.7896
$ws.Range("I86").Value = 1500.4667
$ws.Range("K86").Value = 1500.4667
$ws.Range("M86").Value = -377.4666999999999
$ws.Range("H89").Value = 2162.7896
$ws.Range("I89").Value = 1500.4667
$ws.Range("K89").Value = 7502.3335
$ws.Range("M89").Value = -1886.3335
$ws.Range("H94").Value = 975.5
$ws.Range("I94").Value = 833.5625
$ws.Range("K94").Value = 833.5625
$ws.Range("M94").Value = -382.5625
$ws.Range("H105").Value = 3858.8572
$ws.Range("I105").Value = 3481.25
$ws.Range("K105").Value = 3481.25
$ws.Range("M105").Value = -1734.25
$ws.Range("H134").Value = 908.34485
$ws.Range("I134").Value = 908.34485
$ws.Range("K134").Value = 2725.03455
$ws.Range("M134").Value = -190.0345499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1846
$ws.Range("I31").Value = 1384.3334
$ws.Range("K31").Value = 1384.3334
$ws.Range("M31").Value = -1089.3334
$ws.Range("H34").Value = 1846
$ws.Range("I34").Value = 1384.3334
$ws.Range("K34").Value = 1384.3334
$ws.Range("M34").Value = -1182.3334
$ws.Range("H68").Value = 24996.25
$ws.Range("J68").Value = 24996.25
$ws.Range("L68").Value = 24996.25
$ws.Range("N68").Value = -26494.25
$ws.Range("H71").Value = 24996.25
$ws.Range("J71").Value = 24996.25
$ws.Range("L71").Value = 74988.75
$ws.Range("N71").Value = -82476.75
$ws.Range("H132").Value = 1971.8823
$ws.Range("I132").Value = 1776.2903
$ws.Range("K132").Value = 5328.8709
$ws.Range("M132").Value = -2798.8709

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 86
$ws.Range("J12").Value = 131.85715
$ws.Range("L12").Value = 395.57145
$ws.Range("N12").Value = -741.5714499999999
$ws.Range("H23").Value = 208.44444
$ws.Range("I23").Value = 90
$ws.Range("J23").Value = 242.28572
$ws.Range("K23").Value = 270
$ws.Range("L23").Value = 726.85716
$ws.Range("M23").Value = -35
$ws.Range("N23").Value = -1196.85716
$ws.Range("H47").Value = 551
$ws.Range("I47").Value = 551
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 1653
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -1222
$ws.Range("N47").ClearContents()
$ws.Range("H63").Value = 8499.75
$ws.Range("I63").Value = 3999
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 11997
$ws.Range("L63").Value = 30000
$ws.Range("M63").Value = -11248
$ws.Range("N63").Value = -31498
$ws.Range("H66").Value = 8499.75
$ws.Range("I66").Value = 3999
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 35991
$ws.Range("L66").Value = 90000
$ws.Range("M66").Value = -32247
$ws.Range("N66").Value = -97488
$ws.Range("H86").Value = 2950
$ws.Range("J86").Value = 5000
$ws.Range("L86").Value = 15000
$ws.Range("N86").Value = -17372
$ws.Range("H89").Value = 2950
$ws.Range("J89").Value = 5000
$ws.Range("L89").Value = 45000
$ws.Range("N89").Value = -56856
$ws.Range("H103").Value = 66670172
$ws.Range("I103").Value = 3173.6667
$ws.Range("K103").Value = 9521.000100000001
$ws.Range("M103").Value = -8642.000100000001
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()
$ws.Range("H114").Value = 12500790
$ws.Range("I114").Value = 28572184
$ws.Range("J114").Value = 815.44446
$ws.Range("K114").Value = 85716552
$ws.Range("L114").Value = 2446.33338
$ws.Range("M114").Value = -85713298
$ws.Range("N114").Value = -8954.33338
$ws.Range("H134").Value = 6831.4116
$ws.Range("I134").Value = 2376
$ws.Range("J134").Value = 14999.667
$ws.Range("K134").Value = 7128
$ws.Range("L134").Value = 44999.001
$ws.Range("M134").Value = -2058
$ws.Range("N134").Value = -55139.001
$ws.Range("H141").Value = 3854.9092
$ws.Range("I141").Value = 2906.7
$ws.Range("K141").Value = 8720.099999999999
$ws.Range("M141").Value = -3540.099999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 25216.285
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 25216.285
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H132").Value = 486553.2
$ws.Range("I132").Value = 716087.9
$ws.Range("K132").Value = 2148263.7
$ws.Range("M132").Value = -2145733.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 15000
$ws.Range("J69").Value = 15000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16498
$ws.Range("H72").Value = 15000
$ws.Range("J72").Value = 15000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -52488
$ws.Range("H74").Value = 27694
$ws.Range("J74").Value = 28258.666
$ws.Range("L74").Value = 28258.666
$ws.Range("N74").Value = -30130.666
$ws.Range("H77").Value = 27694
$ws.Range("J77").Value = 28258.666
$ws.Range("L77").Value = 84775.99800000001
$ws.Range("N77").Value = -94135.99800000001
$ws.Range("H81").Value = 10949.429
$ws.Range("J81").Value = 24598.75
$ws.Range("L81").Value = 49197.5
$ws.Range("N81").Value = -51319.5
$ws.Range("H84").Value = 10949.429
$ws.Range("J84").Value = 24598.75
$ws.Range("L84").Value = 245987.5
$ws.Range("N84").Value = -256595.5
$ws.Range("H105").Value = 35474
$ws.Range("J105").Value = 35474
$ws.Range("L105").Value = 35474
$ws.Range("N105").Value = -42462

Write-Output "applied 209 cell changes"